$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (the "statut" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column A holds the status symbol ("⬛" -> "📘"), column B holds its label ("noir" -> "bleu")
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value() -eq "⬛") {
        $cellA.Value = "📘"
    }
    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value() -eq "noir") {
        $cellB.Value = "bleu"
    }
}
